$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 659.625
$ws.Range("I11").Value = 659.625
$ws.Range("K11").Value = 659.625
$ws.Range("M11").Value = -519.625

$ws.Range("H51").Value = 3225.05
$ws.Range("I51").Value = 2750.0625
$ws.Range("K51").Value = 2750.0625
$ws.Range("M51").Value = -2266.0625

$ws.Range("H76").Value = 5078
$ws.Range("I76").Value = 3465
$ws.Range("J76").Value = 7497.5
$ws.Range("K76").Value = 3465
$ws.Range("L76").Value = 7497.5
$ws.Range("M76").Value = -3150
$ws.Range("N76").Value = -8127.5

$ws.Range("H79").Value = 5078
$ws.Range("I79").Value = 3465
$ws.Range("J79").Value = 7497.5
$ws.Range("K79").Value = 3465
$ws.Range("L79").Value = 7497.5
$ws.Range("M79").Value = -2373
$ws.Range("N79").Value = -9681.5

$ws.Range("H88").Value = 2301
$ws.Range("I88").Value = 1850
$ws.Range("J88").Value = 2752
$ws.Range("K88").Value = 1850
$ws.Range("L88").Value = 2752
$ws.Range("M88").Value = -1444
$ws.Range("N88").Value = -3564

$ws.Range("H91").Value = 2301
$ws.Range("I91").Value = 1850
$ws.Range("J91").Value = 2752
$ws.Range("K91").Value = 1850
$ws.Range("L91").Value = 2752
$ws.Range("M91").Value = -446
$ws.Range("N91").Value = -5560

$ws.Range("H99").Value = 1012.17645
$ws.Range("I99").Value = 300.76923
$ws.Range("J99").Value = 3324.25
$ws.Range("K99").Value = 902.30769
$ws.Range("L99").Value = 9972.75
$ws.Range("M99").Value = 595.69231
$ws.Range("N99").Value = -12968.75

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H132").Value = 44946.914
$ws.Range("I132").Value = 28348.324
$ws.Range("K132").Value = 85044.97200000001
$ws.Range("M132").Value = -82514.97200000001

$ws.Range("H133").Value = 65500
$ws.Range("J133").Value = 65500
$ws.Range("L133").Value = 65500
$ws.Range("N133").Value = -75620

$ws.Range("H137").Value = 781353.75
$ws.Range("I137").Value = 19167.133
$ws.Range("J137").Value = 4047867.8
$ws.Range("K137").Value = 57501.399
$ws.Range("L137").Value = 12143603.4
$ws.Range("M137").Value = -54951.399
$ws.Range("N137").Value = -12148703.4

$ws.Range("H138").Value = 4261.0283
$ws.Range("I138").Value = 1976.4445
$ws.Range("J138").Value = 4592.661
$ws.Range("K138").Value = 5929.333500000001
$ws.Range("L138").Value = 13777.983
$ws.Range("M138").Value = -789.3335000000006
$ws.Range("N138").Value = -24057.983

$ws.Range("H140").Value = 70543.63
$ws.Range("J140").Value = 70543.63
$ws.Range("L140").Value = 70543.63
$ws.Range("N140").Value = -80903.63

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2353.5
$ws.Range("I2").Value = 2854.1667
$ws.Range("J2").Value = 1602.5
$ws.Range("K2").Value = 2854.1667
$ws.Range("L2").Value = 1602.5
$ws.Range("M2").Value = -2741.1667
$ws.Range("N2").Value = -1828.5

$ws.Range("H32").Value = 3481.4878
$ws.Range("I32").Value = 3189.9412
$ws.Range("K32").Value = 3189.9412
$ws.Range("M32").Value = -2902.9412

$ws.Range("H61").Value = 2636.182
$ws.Range("I61").Value = 2299.8
$ws.Range("K61").Value = 2299.8
$ws.Range("M61").Value = -2087.8

$ws.Range("H97").Value = 3639.2144
$ws.Range("I97").Value = 3565
$ws.Range("J97").Value = 3824.75
$ws.Range("K97").Value = 3565
$ws.Range("L97").Value = 3824.75
$ws.Range("M97").Value = -3069
$ws.Range("N97").Value = -4816.75

$ws.Range("H116").Value = 2353.5
$ws.Range("I116").Value = 2854.1667
$ws.Range("J116").Value = 1602.5
$ws.Range("K116").Value = 2854.1667
$ws.Range("L116").Value = 1602.5
$ws.Range("M116").Value = -560.1667000000002
$ws.Range("N116").Value = -6190.5

$ws.Range("H136").Value = 2636.182
$ws.Range("I136").Value = 2299.8
$ws.Range("K136").Value = 6899.400000000001
$ws.Range("M136").Value = -4349.400000000001

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2353.5
$ws.Range("I3").Value = 2854.1667
$ws.Range("J3").Value = 1602.5
$ws.Range("K3").Value = 2854.1667
$ws.Range("L3").Value = 1602.5
$ws.Range("M3").Value = -2740.1667
$ws.Range("N3").Value = -1830.5

$ws.Range("H134").Value = 1814.4333
$ws.Range("I134").Value = 1590.4642
$ws.Range("K134").Value = 4771.392599999999
$ws.Range("M134").Value = -2236.392599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 12482.714
$ws.Range("I32").Value = 10338.333
$ws.Range("K32").Value = 10338.333
$ws.Range("M32").Value = -10022.333

$ws.Range("H58").Value = 3247.138
$ws.Range("I58").Value = 2367.8
$ws.Range("K58").Value = 2367.8
$ws.Range("M58").Value = -2164.8

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H118").Value = 203999.5
$ws.Range("J118").Value = 203999.5
$ws.Range("L118").Value = 203999.5
$ws.Range("N118").Value = -207313.5

$ws.Range("H132").Value = 1907.625
$ws.Range("I132").Value = 1743.3549
$ws.Range("K132").Value = 5230.0647
$ws.Range("M132").Value = -2700.0647

$ws.Range("H134").Value = 21148.875
$ws.Range("I134").Value = 23456.035
$ws.Range("K134").Value = 70368.105
$ws.Range("M134").Value = -67833.105

$ws.Range("H136").Value = 3247.138
$ws.Range("I136").Value = 2367.8
$ws.Range("K136").Value = 7103.400000000001
$ws.Range("M136").Value = -4553.400000000001

$ws.Range("H140").Value = 79667
$ws.Range("J140").Value = 79667
$ws.Range("L140").Value = 79667
$ws.Range("N140").Value = -90027

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 62501948
$ws.Range("I18").Value = 62501948
$ws.Range("K18").Value = 187505844
$ws.Range("M18").Value = -187505675

$ws.Range("H140").Value = 4085.0952
$ws.Range("J140").Value = 3992.7693
$ws.Range("L140").Value = 11978.3079
$ws.Range("N140").Value = -22338.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6249.75
$ws.Range("I102").Value = 6249.75
$ws.Range("K102").Value = 6249.75
$ws.Range("M102").Value = -4627.75

$ws.Range("H109").Value = 38664.668
$ws.Range("J109").Value = 38664.668
$ws.Range("L109").Value = 38664.668
$ws.Range("N109").Value = -40744.668

$ws.Range("H132").Value = 19583.906
$ws.Range("I132").Value = 21887.785
$ws.Range("K132").Value = 65663.355
$ws.Range("M132").Value = -63133.355

$ws.Range("H135").Value = 48562.25
$ws.Range("J135").Value = 48562.25
$ws.Range("L135").Value = 48562.25
$ws.Range("N135").Value = -58702.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3940.2222
$ws.Range("I16").Value = 4110.846
$ws.Range("J16").Value = 3496.6
$ws.Range("K16").Value = 4110.846
$ws.Range("L16").Value = 3496.6
$ws.Range("M16").Value = -3940.846
$ws.Range("N16").Value = -3836.6

$ws.Range("H55").Value = 333.55554
$ws.Range("I55").Value = 218.6
$ws.Range("J55").Value = 477.25
$ws.Range("K55").Value = 218.6
$ws.Range("L55").Value = 477.25
$ws.Range("M55").Value = -45.59999999999999
$ws.Range("N55").Value = -823.25

$ws.Range("H138").Value = 78500
$ws.Range("J138").Value = 78500
$ws.Range("L138").Value = 78500
$ws.Range("N138").Value = -88780

$ws.Range("H141").Value = 80660.336
$ws.Range("J141").Value = 80660.336
$ws.Range("L141").Value = 80660.336
$ws.Range("N141").Value = -91020.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10333.667

$ws.Range("H63").Value = 27599.4
$ws.Range("J63").Value = 29499.75
$ws.Range("L63").Value = 29499.75
$ws.Range("N63").Value = -30747.75

$ws.Range("H66").Value = 27599.4
$ws.Range("J66").Value = 29499.75
$ws.Range("L66").Value = 88499.25
$ws.Range("N66").Value = -94739.25

$ws.Range("H122").Value = 4586.8
$ws.Range("I122").Value = 2550.8696
$ws.Range("K122").Value = 7652.6088
$ws.Range("M122").Value = -5202.6088

$ws.Range("H136").Value = 1488.16
$ws.Range("I136").Value = 1191.1364
$ws.Range("K136").Value = 3573.4092
$ws.Range("M136").Value = -1023.4092

$ws.Range("H138").Value = 84524.5
$ws.Range("J138").Value = 79050
$ws.Range("L138").Value = 79050
$ws.Range("N138").Value = -89330

$ws.Range("H139").Value = 72316.25
$ws.Range("J139").Value = 72316.25
$ws.Range("L139").Value = 72316.25
$ws.Range("N139").Value = -82596.25
